# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the "30e077c1-baa6-4a0a-818c-669dbbb83c18" entry:
#   - Status changes from "Ready for handoff" to "Handback transform failed"
#     (Overview!B3, Overview!C3, zh-cn!C3, de-de!C3)
#   - A new "Error Detail" message is recorded in column K of row 3 on the
#     zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K3").Value = "Handback file name: vggu2zue.xu3 is different with handoff file name: 30e077c1-baa6-4a0a-818c-669dbbb83c18.4efc8708f9f31096b1556d6f0a6d5c2fc633cc7e.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K3").Value = "Handback file name: vggu2zue.xu3 is different with handoff file name: 30e077c1-baa6-4a0a-818c-669dbbb83c18.4efc8708f9f31096b1556d6f0a6d5c2fc633cc7e.de-de."
